$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVTs")

# ---- Row 2 ----
$ws.Range("B2").Value = 'Display text string'
$ws.Range("C2").Value = 'Check out for numeric value'
$ws.Range("D2").Value = "1. Drag 'Values' column from 'Measure data' table into 'Field' field"
$ws.Range("E2").Value = 'Numeric value should be displayed'

# ---- Row 3 ----
$ws.Range("C3").Value = 'Check output for text value'
$ws.Range("D3").Value = "1. Drag 'Column1' column from 'Column data' table into 'Field' field of text wrapper visual`n2. Select slicer visual from 'Visualizations' pane. Drag 'Column1' column from 'Column data' table into 'Field' field`n3. Click on 'What is your salary range?' checkbox"
$ws.Range("E3").Value = "What is your salary range?' text should be displayed in the visual"

# ---- Row 4 ----
$ws.Range("C4").Value = 'Display error message part 1'
$ws.Range("D4").Value = "1. Drag 'Column1' column from 'Column data' table into 'Field' field of text wrapper visual`n2. Select slicer visual from 'Visualizations' pane. Drag 'Column1' column from 'Column data' table into 'Field' field"
$ws.Range("E4").Value = '"Query returned more than one row, please filter data to return one row" text should be displayed in the visual'

# ---- Row 5 ----
$ws.Range("C5").Value = 'Display error message part 2'
$ws.Range("D5").Value = "1. Drag 'Column1' column from 'Column data' table into 'Field' field of text wrapper visual`n2. Select slicer visual from 'Visualizations' pane. Drag 'Column1' column from 'Column data' table into 'Field' field`n3. Click on 'What is your salary range?' checkbox and 'What is your job title' checkbox"
$ws.Range("E5").Value = '"Query returned more than one row, please filter data to return one row" text should be displayed in the visual'

# ---- Row 6 (existing row, content replaced) ----
$ws.Range("B6").Value = 'Text settings'
$ws.Range("C6").Value = 'Update text color, text size'
$ws.Range("D6").Value = "1. Go to formatting pane`n2. Go to 'Text settings' option`n3. Update color to 'blue' `n4. Update text size to '25'"
$ws.Range("E6").Value = "1. Text color will be set to 'blue'`n2. Text size will be set to 25 pt."
$ws.Rows.Item(6).RowHeight = 60

# ---- Row 7 (new row) ----
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 'Static Text'
$ws.Range("C7").Value = 'Features for static text '
$ws.Range("D7").Value = "1. Go to formattig pane`n2. Go to 'Static Text Settigs' option`n3. Check Font style`n4. Check Font family`n5. Check Text Highlighter`n6. Check Colon`n7. Check Position of the static text(Text Position)- suffix or prefix"
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value = "1. Bold, Italics and Underline are by default off.`n2. Colon will be set to 'off' by default. Segoe UI Semibold is the default font family.`n3. Text Position will be set to suffix by default. Text highlighter default value is set to white. Add any text to ""Text to append"" field say 'Hello world', the text 'Hello world' will be appended as suffix with a colon by default."
$ws.Range("E7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 120

# ---- Row 8 (new row) ----
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 'Dynamic Text'
$ws.Range("C8").Value = 'User will drag the column to the Field'
$ws.Range("D8").Value = "1. Go to Formatting pane`n2. Go to 'Dynamic Text Settings' option`n3.Check Font style `n4. Check Font Family`n5. Check Text Highlighter"
$ws.Range("D8").WrapText = $true
$ws.Range("E8").Value = "1. Bold, Italics and Underline are by default off.`n2. Segoe UI Semibold is the default font family.`n3. Text highlighter default value is set to white."
$ws.Rows.Item(8).RowHeight = 71.25

# ---- Row 9 (new row) ----
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = 'Colon'
$ws.Range("B9").WrapText = $true
$ws.Range("C9").Value = 'toggle option is there'
$ws.Range("D9").Value = "1. Go to formattig pane`n2. Go to 'Static Text Settigs' option`n3. Switch toggle of 'Show Colon' to 'Off'"
$ws.Range("D9").WrapText = $true
$ws.Range("E9").Value = 'Colon appears with the appearance of static text and then you can remove the colon with the toggle button'
$ws.Range("E9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 45

# ---- Checklist sheet selection (select first, before BVTs so BVTs stays the active tab) ----
$ws2 = $wb.Worksheets.Item("Checklist")
$ws2.Range("C21:C25").Select()

# ---- BVTs sheet: final active selection ----
$ws.Activate()
$ws.Range("D9").Select()
